$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New compiled results data: Target Dimension (A) and M1 (B)
$data = @(
    @(1, 171.1777899565496),
    @(2, 95.02867125292212),
    @(3, 94.88894051753742),
    @(5, 122.4481994433494),
    @(6, 182.7437561640268),
    @(7, 190.6306695265143),
    @(8, 256.7669239621679),
    @(10, 283.0170036981232)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}
